$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oldest day (2025-10-11) -- row 2 -- which shifts every
# subsequent row up by one and keeps all the Non-HTTPS/HTTPS numbers
# correctly aligned with their (now shifted) dates.
$ws.Rows.Item(2).Delete()

# Append the two newest days to the bottom of the table. Column A must
# stay plain text (not get auto-converted to a date serial), so mark the
# cells as Text before typing the date strings into them.
$ws.Cells.Item(90,1).NumberFormat = "@"
$ws.Cells.Item(90,1).Value = "2026-01-08"
$ws.Cells.Item(90,2).Value = 0
$ws.Cells.Item(90,3).Value = 27

$ws.Cells.Item(91,1).NumberFormat = "@"
$ws.Cells.Item(91,1).Value = "2026-01-09"
$ws.Cells.Item(91,2).Value = 0
$ws.Cells.Item(91,3).Value = 27
